# Apply the daily cryptos data refresh to Sheet1.
# This mirrors the automated "Updated cryptos list ... with GitHub Actions" commit:
# most rows keep their Coin/Link but get refreshed Price / Volume(1h) figures;
# two row-pairs (36/37 and 49/50) were re-ranked, swapping which coin occupies
# which row (along with its own refreshed price/volume).
#
# The Price column stores values as plain text (e.g. "318.70", "0.0932",
# "45.342.09") even though many of them look like numbers. If we just assign
# a plain numeric-looking string to .Value, Excel auto-converts it to a
# floating point number (losing trailing zeros / exact text / introducing
# binary rounding noise). To keep these as text, exactly like the source
# file, we force the cell's number format to Text ("@") before assigning the
# value, then restore the cell style afterwards so no stray formatting is
# left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "45.342.09"
$ws.Range("E2").Value = "  +7.04%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.389.08"
$ws.Range("E3").Value = "  +4.92%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.44%  "

# Row 5 - Solana
Set-TextValue "D5" "112.38"
$ws.Range("E5").Value = "  +9.09%  "

# Row 6 - BNB
Set-TextValue "D6" "318.70"
$ws.Range("E6").Value = "  +2.98%  "

# Row 7 - XRP
Set-TextValue "D7" "0.636"
$ws.Range("E7").Value = "  +2.75%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.09%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +5.46%  "

# Row 10 - Avalanche
$ws.Range("E10").Value = "  +10.21%  "

# Row 11 - Dogecoin
Set-TextValue "D11" "0.0932"
$ws.Range("E11").Value = "  +4.10%  "

# Row 12 - Polkadot
Set-TextValue "D12" "8.72"
$ws.Range("E12").Value = "  +6.63%  "

# Row 13 - Polygon
$ws.Range("E13").Value = "  +5.74%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  +1.91%  "

# Row 15 - Chainlink
Set-TextValue "D15" "15.90"
$ws.Range("E15").Value = "  +6.74%  "

# Row 16 - WrappedliquidstakedEther2.0
Set-TextValue "D16" "2.747.34"
$ws.Range("E16").Value = "  +4.84%  "

# Row 17 - WrappedEther
Set-TextValue "D17" "2.391.25"
$ws.Range("E17").Value = "  +5.07%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "45.337.20"
$ws.Range("E18").Value = "  +7.17%  "

# Row 19 - Uniswap
Set-TextValue "D19" "7.70"
$ws.Range("E19").Value = "  +6.79%  "

# Row 20 - ShibaInu
$ws.Range("E20").Value = "  +4.31%  "

# Row 21 - InternetComputer(DFINITY)
Set-TextValue "D21" "13.18"
$ws.Range("E21").Value = "  +2.36%  "

# Row 22 - Litecoin
$ws.Range("E22").Value = "  +3.54%  "

# Row 23 - PancakeSwap
Set-TextValue "D23" "3.53"
$ws.Range("E23").Value = "  +5.33%  "

# Row 24 - BitcoinCash
Set-TextValue "D24" "269.58"
$ws.Range("E24").Value = "  +3.33%  "

# Row 25 - ImmutableX
Set-TextValue "D25" "2.34"
$ws.Range("E25").Value = "  +8.46%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  -0.78%  "

# Row 27 - Cosmos
Set-TextValue "D27" "11.31"

# Row 28 - Filecoin
Set-TextValue "D28" "7.54"
$ws.Range("E28").Value = "  +11.16%  "

# Row 29 - Toncoin
$ws.Range("E29").Value = "  -4.03%  "

# Row 30 - EthereumClassic
Set-TextValue "D30" "22.93"

# Row 31 - InjectiveProtocol
Set-TextValue "D31" "38.94"
$ws.Range("E31").Value = "  +9.18%  "

# Row 32 - Hedera
Set-TextValue "D32" "0.0947"
$ws.Range("E32").Value = "  +11.53%  "

# Row 33 - Monero
Set-TextValue "D33" "170.06"
$ws.Range("E33").Value = "  +3.78%  "

# Row 34 - WEMIXToken
$ws.Range("E34").Value = "  +17.27%  "

# Row 35 - Stellar
$ws.Range("E35").Value = "  +3.66%  "

# Row 36 - was Kaspa, now RenderToken (re-ranked)
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D36" "4.94"
$ws.Range("E36").Value = "  +11.34%  "

# Row 37 - was RenderToken, now Kaspa (re-ranked)
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D37" "0.118"
$ws.Range("E37").Value = "  +7.09%  "

# Row 38 - LidoDAOToken
Set-TextValue "D38" "3.11"
$ws.Range("E38").Value = "  +15.01%  "

# Row 39 - VeChain
$ws.Range("E39").Value = "  +6.21%  "

# Row 40 - NEARProtocol
Set-TextValue "D40" "3.97"
$ws.Range("E40").Value = "  +8.52%  "

# Row 41 - ARBITRUM
$ws.Range("E41").Value = "  +13.43%  "

# Row 42 - BitcoinSV
Set-TextValue "D42" "105.15"
$ws.Range("E42").Value = "  +5.95%  "

# Row 43 - Algorand
Set-TextValue "D43" "0.241"
$ws.Range("E43").Value = "  +7.45%  "

# Row 44 - Celestia
Set-TextValue "D44" "13.56"
$ws.Range("E44").Value = "  +14.44%  "

# Row 45 - MultiversX
Set-TextValue "D45" "71.32"
$ws.Range("E45").Value = "  +4.70%  "

# Row 46 - FirstDigitalUSD
$ws.Range("E46").Value = "  +0.19%  "

# Row 47 - Aave
Set-TextValue "D47" "117.87"
$ws.Range("E47").Value = "  +7.83%  "

# Row 48 - THORChain
Set-TextValue "D48" "5.84"
$ws.Range("E48").Value = "  +14.69%  "

# Row 49 - was FraxShare, now MinaProtocolToken (re-ranked)
$ws.Range("B49").Value = "MinaProtocolToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/J7st_qGwz+minaprotocoltoken-mina"
Set-TextValue "D49" "1.64"
$ws.Range("E49").Value = "  +20.55%  "

# Row 50 - was MinaProtocolToken, now FraxShare (re-ranked)
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D50" "9.40"
$ws.Range("E50").Value = "  +9.63%  "

# Row 51 - ordi
Set-TextValue "D51" "78.90"
$ws.Range("E51").Value = "  +3.89%  "
